# Apply crypto price/volume updates from the Thu Oct 26 15:45:52 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as literal text (matches the source inlineStr cells), even when the
# new value looks like a number (e.g. "220.57") or has multiple dots (e.g. "33.950.03").
# Flipping to NumberFormat "@" before the write defeats Excel's automatic number/date
# detection, and ClearFormats() afterwards restores the original (default) cell style so
# no stray number-format style is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "33.950.03"
Set-TextValue $ws.Range("E2") "  -2.26%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.774.09"
Set-TextValue $ws.Range("E3") "  -1.62%  "

# Row 4
Set-TextValue $ws.Range("E4") "  +0.24%  "

# Row 5
Set-TextValue $ws.Range("D5") "220.57"
Set-TextValue $ws.Range("E5") "  -2.49%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -2.31%  "

# Row 7
Set-TextValue $ws.Range("E7") "  +0.20%  "

# Row 8
Set-TextValue $ws.Range("D8") "31.13"
Set-TextValue $ws.Range("E8") "  -5.82%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.284"
Set-TextValue $ws.Range("E9") "  -0.66%  "

# Row 10
Set-TextValue $ws.Range("E10") "  +4.36%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0920"
Set-TextValue $ws.Range("E11") "  -1.72%  "

# Row 12
Set-TextValue $ws.Range("D12") "2.032.31"
Set-TextValue $ws.Range("E12") "  -1.51%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.781.25"
Set-TextValue $ws.Range("E13") "  -1.21%  "

# Row 14
Set-TextValue $ws.Range("D14") "10.47"
Set-TextValue $ws.Range("E14") "  -7.83%  "

# Row 15
Set-TextValue $ws.Range("B15") "WrappedBTC"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D15") "33.949.02"
Set-TextValue $ws.Range("E15") "  -2.30%  "

# Row 16
Set-TextValue $ws.Range("B16") "Polygon"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D16") "0.618"
Set-TextValue $ws.Range("E16") "  -3.95%  "

# Row 17
Set-TextValue $ws.Range("D17") "4.18"
Set-TextValue $ws.Range("E17") "  -2.64%  "

# Row 18
Set-TextValue $ws.Range("D18") "67.60"
Set-TextValue $ws.Range("E18") "  -3.01%  "

# Row 19
Set-TextValue $ws.Range("D19") "242.51"
Set-TextValue $ws.Range("E19") "  -5.60%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -0.29%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.999"
Set-TextValue $ws.Range("E21") "  +0.35%  "

# Row 22
Set-TextValue $ws.Range("D22") "10.47"
Set-TextValue $ws.Range("E22") "  -0.50%  "

# Row 23
Set-TextValue $ws.Range("E23") "  -5.39%  "

# Row 24
Set-TextValue $ws.Range("E24") "  -1.48%  "

# Row 25
Set-TextValue $ws.Range("D25") "157.49"
Set-TextValue $ws.Range("E25") "  -0.89%  "

# Row 26
Set-TextValue $ws.Range("D26") "16.24"
Set-TextValue $ws.Range("E26") "  -1.93%  "

# Row 27
Set-TextValue $ws.Range("D27") "6.95"
Set-TextValue $ws.Range("E27") "  -3.14%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -3.21%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.00"
Set-TextValue $ws.Range("E29") "  +0.20%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.0514"
Set-TextValue $ws.Range("E30") "  -1.65%  "

# Row 31
Set-TextValue $ws.Range("B31") "PancakeSwap"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.20"
Set-TextValue $ws.Range("E31") "  -0.32%  "

# Row 32
Set-TextValue $ws.Range("B32") "Filecoin"
Set-TextValue $ws.Range("C32") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "3.66"
Set-TextValue $ws.Range("E32") "  -4.41%  "

# Row 33
Set-TextValue $ws.Range("D33") "3.49"
Set-TextValue $ws.Range("E33") "  -3.62%  "

# Row 34
Set-TextValue $ws.Range("D34") "1.82"
Set-TextValue $ws.Range("E34") "  -5.10%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.391.66"
Set-TextValue $ws.Range("E35") "  -4.77%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.623"
Set-TextValue $ws.Range("E37") "  -2.38%  "

# Row 38
Set-TextValue $ws.Range("E38") "  -3.24%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.926"
Set-TextValue $ws.Range("E39") "  +2.09%  "

# Row 40
Set-TextValue $ws.Range("E40") "  +0.13%  "

# Row 41
Set-TextValue $ws.Range("D41") "78.57"
Set-TextValue $ws.Range("E41") "  -6.19%  "

# Row 42
Set-TextValue $ws.Range("E42") "  -5.62%  "

# Row 43
Set-TextValue $ws.Range("E43") "  -0.97%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.0489"
Set-TextValue $ws.Range("E44") "  -3.80%  "

# Row 45
Set-TextValue $ws.Range("D45") "5.81"
Set-TextValue $ws.Range("E45") "  -4.16%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.03"
Set-TextValue $ws.Range("E46") "  -1.02%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.929.64"
Set-TextValue $ws.Range("E47") "  -1.72%  "

# Row 48
Set-TextValue $ws.Range("D48") "103.94"
Set-TextValue $ws.Range("E48") "  +1.16%  "

# Row 49
Set-TextValue $ws.Range("E49") "  -0.24%  "

# Row 50
Set-TextValue $ws.Range("D50") "11.68"
Set-TextValue $ws.Range("E50") "  -3.42%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.0₆0120"
Set-TextValue $ws.Range("E51") "  -1.40%  "
